$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" target-cluster rows were dropped from the NATMI LR-pair
# table (one per sending cluster: FAPs->Resolving-Mac was row 6, MuSCs->Resolving-Mac
# was row 11). Delete bottom row first so the row-6 index stays valid.
$ws.Rows("11").Delete()
$ws.Rows("6").Delete()

# Refresh the remaining rows (2-9) with the re-run (new TPM) NATMI values.

# Row 2
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Btc"
$ws.Range("C2").Value2 = "Erbb2"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.4631816666666667
$ws.Range("H2").Value2 = 1.389545
$ws.Range("I2").Value2 = 0.3800727954645477
$ws.Range("J2").Value2 = 0.47906870267432
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.271905
$ws.Range("N2").Value2 = 6.815715
$ws.Range("O2").Value2 = 0.2806394474136332
$ws.Range("P2").Value2 = 0.3070885820898414
$ws.Range("Q2").Value2 = 1.052304744408333
$ws.Range("R2").Value2 = 9.470742699675
$ws.Range("S2").Value2 = 0.1066634192961255
$ws.Range("T2").Value2 = 0.1471165286278767

# Row 3
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Btc"
$ws.Range("C3").Value2 = "Erbb2"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.4631816666666667
$ws.Range("H3").Value2 = 1.389545
$ws.Range("I3").Value2 = 0.3800727954645477
$ws.Range("J3").Value2 = 0.47906870267432
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 3.687724
$ws.Range("N3").Value2 = 11.063172
$ws.Range("O3").Value2 = 0.4555299739971492
$ws.Range("P3").Value2 = 0.4984618345831706
$ws.Range("Q3").Value2 = 1.708086148526667
$ws.Range("R3").Value2 = 15.37277533674
$ws.Range("S3").Value2 = 0.1731345506349892
$ws.Range("T3").Value2 = 0.238797464426421

# Row 4
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Btc"
$ws.Range("C4").Value2 = "Erbb2"
$ws.Range("D4").Value2 = "Inflammatory-Mac"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.4631816666666667
$ws.Range("H4").Value2 = 1.389545
$ws.Range("I4").Value2 = 0.3800727954645477
$ws.Range("J4").Value2 = 0.47906870267432
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.044076
$ws.Range("N4").Value2 = 0.132228
$ws.Range("O4").Value2 = 0.005444534117493161
$ws.Range("P4").Value2 = 0.005957659472641616
$ws.Range("Q4").Value2 = 0.02041519514
$ws.Range("R4").Value2 = 0.18373675626
$ws.Range("S4").Value2 = 0.00206931930203773
$ws.Range("T4").Value2 = 0.002854128194533792

# Row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Btc"
$ws.Range("C5").Value2 = "Erbb2"
$ws.Range("D5").Value2 = "MuSCs"
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.4631816666666667
$ws.Range("H5").Value2 = 1.389545
$ws.Range("I5").Value2 = 0.3800727954645477
$ws.Range("J5").Value2 = 0.47906870267432
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 2.0917535
$ws.Range("N5").Value2 = 4.183507000000001
$ws.Range("O5").Value2 = 0.2583860444717245
$ws.Range("P5").Value2 = 0.1884919238543463
$ws.Range("Q5").Value2 = 0.9688618723858334
$ws.Range("R5").Value2 = 5.813171234315001
$ws.Range("S5").Value2 = 0.09820550623139528
$ws.Range("T5").Value2 = 0.09030058142548841

# Row 6
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Btc"
$ws.Range("C6").Value2 = "Erbb2"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.5
$ws.Range("G6").Value2 = 0.755484
$ws.Range("H6").Value2 = 1.510968
$ws.Range("I6").Value2 = 0.6199272045354524
$ws.Range("J6").Value2 = 0.52093129732568
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 2.271905
$ws.Range("N6").Value2 = 6.815715
$ws.Range("O6").Value2 = 0.2806394474136332
$ws.Range("P6").Value2 = 0.3070885820898414
$ws.Range("Q6").Value2 = 1.71638787702
$ws.Range("R6").Value2 = 10.29832726212
$ws.Range("S6").Value2 = 0.1739760281175077
$ws.Range("T6").Value2 = 0.1599720534619646

# Row 7
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Btc"
$ws.Range("C7").Value2 = "Erbb2"
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.5
$ws.Range("G7").Value2 = 0.755484
$ws.Range("H7").Value2 = 1.510968
$ws.Range("I7").Value2 = 0.6199272045354524
$ws.Range("J7").Value2 = 0.52093129732568
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 3.687724
$ws.Range("N7").Value2 = 11.063172
$ws.Range("O7").Value2 = 0.4555299739971492
$ws.Range("P7").Value2 = 0.4984618345831706
$ws.Range("Q7").Value2 = 2.786016478416
$ws.Range("R7").Value2 = 16.716098870496
$ws.Range("S7").Value2 = 0.28239542336216
$ws.Range("T7").Value2 = 0.2596643701567495

# Row 8
$ws.Range("A8").Value2 = "MuSCs"
$ws.Range("B8").Value2 = "Btc"
$ws.Range("C8").Value2 = "Erbb2"
$ws.Range("D8").Value2 = "Inflammatory-Mac"
$ws.Range("E8").Value2 = 1
$ws.Range("F8").Value2 = 0.5
$ws.Range("G8").Value2 = 0.755484
$ws.Range("H8").Value2 = 1.510968
$ws.Range("I8").Value2 = 0.6199272045354524
$ws.Range("J8").Value2 = 0.52093129732568
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.044076
$ws.Range("N8").Value2 = 0.132228
$ws.Range("O8").Value2 = 0.005444534117493161
$ws.Range("P8").Value2 = 0.005957659472641616
$ws.Range("Q8").Value2 = 0.033298712784
$ws.Range("R8").Value2 = 0.199792276704
$ws.Range("S8").Value2 = 0.003375214815455432
$ws.Range("T8").Value2 = 0.003103531278107823

# Row 9
$ws.Range("A9").Value2 = "MuSCs"
$ws.Range("B9").Value2 = "Btc"
$ws.Range("C9").Value2 = "Erbb2"
$ws.Range("D9").Value2 = "MuSCs"
$ws.Range("E9").Value2 = 1
$ws.Range("F9").Value2 = 0.5
$ws.Range("G9").Value2 = 0.755484
$ws.Range("H9").Value2 = 1.510968
$ws.Range("I9").Value2 = 0.6199272045354524
$ws.Range("J9").Value2 = 0.52093129732568
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 2.0917535
$ws.Range("N9").Value2 = 4.183507000000001
$ws.Range("O9").Value2 = 0.2583860444717245
$ws.Range("P9").Value2 = 0.1884919238543463
$ws.Range("Q9").Value2 = 1.580286301194
$ws.Range("R9").Value2 = 6.321145204776001
$ws.Range("S9").Value2 = 0.1601805382403293
$ws.Range("T9").Value2 = 0.09819134242885794
